# cronograma_funcionalidades.xlsx - "coisas a fazer" sheet reshuffle
# Commit: criando função de enviar email ao logar
#
# The TODO list on the "coisas a fazer" sheet is reordered: the item
# "Colocar a Quant. De crismandos no card de grupos da crisma (pagina
# Grupos)" moves from the top of the list to the bottom, the remaining
# items shift up one row, and a brand-new item #8 is appended:
# "Criar a rota e paginas para editar informações de um grupo".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("coisas a fazer")

# Copy the formatting used by the other numbered rows (A: centered
# number, B: bordered + wrapped text) onto the new row 9 before it gets
# its own values, so it matches rows 2-8's look.
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)

# New text for rows 3-8 (each existing item shifts up one slot; the old
# row 3 item becomes the new row 8 / item 7), plus the brand-new row 9.
$ws.Range("B3").Value = "Criar página de perfi, com as informações do usuário (nele terá a opção de editar infor., foto de perfil (?), senha e email)"
$ws.Range("B4").Value = "Melhorar ou modificar o sistema de login e cadastro de catequista"
$ws.Range("B5").Value = "Página p/ relatar bugs"
$ws.Range("B6").Value = "Mehorar o digitamento das informações (por exemplo: permitir onde é para aceitar numero, só aceite numero e já fique formatado)"
$ws.Range("B7").Value = "Modificar algumas coisas no banco de dados (fazer melhoramentos e adicionar colunas ( adicionar colunas em ""crismandos"" como por exemplo estado civil, se possui filhos, cidade, etc))"
$ws.Range("B8").Value = "Colocar a Quant. De crismandos no card de grupos da crisma (pagina Grupos)"
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Criar a rota e paginas para editar informações de um grupo"

# Row heights follow the text that now occupies each row.
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 57.6
$ws.Rows.Item(8).RowHeight = 28.8

# The old trailing blank placeholder row 10 is no longer needed now that
# row 9 carries real content - wipe it completely (contents + formatting)
# so it drops out of the sheet, leaving rows 11-13 as the remaining
# blank placeholders.
$ws.Range("B10").ClearContents()
$ws.Range("B10").ClearFormats()

# Update the remembered selection to match the edited file.
$ws.Range("E7").Select()
